$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.429.87'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '3.846.13'
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.08'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.67'
$ws.Range("E6").Value = '  -2.72%  '
$ws.Range("D7").Value = '3.845.73'
$ws.Range("E7").Value = '  +2.88%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -2.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.32'
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.81'
$ws.Range("E13").Value = '  -3.09%  '
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").Value = '4.481.32'
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").Value = '3.828.48'
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("D17").Value = '68.627.98'
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.53'
$ws.Range("E18").Value = '  +2.53%  '
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.12'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.26'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.66'
$ws.Range("E22").Value = '  -1.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.717'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000161'
$ws.Range("E24").Value = '  +6.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.92'
$ws.Range("E25").Value = '  -1.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.09'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.93'
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").Value = '3.992.80'
$ws.Range("E31").Value = '  +2.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.83'
$ws.Range("E32").Value = '  -4.24%  '
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '32.02'
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").Value = '3.789.26'
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.87'
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.96'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '430.69'
$ws.Range("E43").Value = '  +2.01%  '
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.36'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.95'
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").Value = '2.828.15'
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.03'
$ws.Range("E50").Value = '  +14.01%  '
$ws.Range("E51").Value = '  +0.69%  '
